# Weekly update: insert two new price records at the top of the data
# (rows 735-736), pushing all existing rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 735, shifting
# all subsequent rows (old 735..826) down to (new 737..828).
$ws.Rows.Item(735).Insert()
$ws.Rows.Item(735).Insert()

# --- New row 735 ---
$ws.Range("A735").Value = 3
$ws.Range("B735").Value = "Femacal de La Calera"
$ws.Range("C735").Value = "Coquimbo"
$ws.Range("D735").Value = 45142
$ws.Range("E735").Value = 5
$ws.Range("F735").Value = 100112003
$ws.Range("G735").Value = "Ajo"
$ws.Range("H735").Value = "Chino"
$ws.Range("I735").Value = "Primera"
$ws.Range("J735").Value = 40
$ws.Range("K735").Value = 18000
$ws.Range("L735").Value = 18000
$ws.Range("M735").Value = 18000
$ws.Range("N735").Value = "`$/caja 10 kilos"
$ws.Range("O735").Value = "China"
$ws.Range("P735").Value = 1800
$ws.Range("Q735").Value = 10
$ws.Range("R735").Value = "Hortaliza"

# --- New row 736 ---
$ws.Range("A736").Value = 3
$ws.Range("B736").Value = "Femacal de La Calera"
$ws.Range("C736").Value = "Coquimbo"
$ws.Range("D736").Value = 45142
$ws.Range("E736").Value = 5
$ws.Range("F736").Value = 100112003
$ws.Range("G736").Value = "Ajo"
$ws.Range("H736").Value = "Chino"
$ws.Range("I736").Value = "Primera"
$ws.Range("J736").Value = 40
$ws.Range("K736").Value = 21000
$ws.Range("L736").Value = 21000
$ws.Range("M736").Value = 21000
$ws.Range("N736").Value = "`$/malla 10 kilos"
$ws.Range("O736").Value = "China"
$ws.Range("P736").Value = 2100
$ws.Range("Q736").Value = 10
$ws.Range("R736").Value = "Hortaliza"
